$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 16 rush/pass yardage logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "2 4 9 2 5 19 5 2 0 6 4 0 1 1 4 15 3 3 2 1 8 21 3 5 7 12 5 12 4 -2 7 4 5 5 2 8 11 4 5 13 75 18 4 11 3 5 3 9 4 4 4 7 1 3 14 3 3 -1 0 12 4 2 1 4 4 6 -2 13 2 -3 13 2 2 9 5 3 3 19 5 2 0 0 9 23 0 8 -5 3 0 3 3 -4 5 2 2 -3 0 6 1 11 3 12 4 5 -1 2 0 2 5 -2 25 0 6 1 0 7 4 1 7 2 -3 8 3 20 4 3 1 3 5 5 5 2 2 0 8 2 4 4 0 9 7 5 6 0 3 -3 4 9 -2 5 13 1 4 6 10 9 4 3 3 7 2 3 8 1 2 6 4 8 3 5 1 2 9 1 7 11 6 2 2 9 2 11 6 3 2 5 3 7 2 3 0 4 4 2 3 8 4 -2 9 0 0 0 2 0 4 2 1 7 0 7 -3 4 0 0 3 6 1 4 3 5 2 5 2 2 1 20 5 4 0 0 7 2 4 1 1 15 4 1 2 4 4 2 9 3 6 1 0 14 1 9 8 1 4 0 3 5 1 10 2 1 7 2 2 5 7 3 11 4 8 2 0 17 4 4 1 13 2 16 1 -4 1 4 10 4 1 5 3 6 1 1 1 3 6 7 3 18 3 3 8 -4 1 3 77 1 5 3 5 3 6 2 5 1 2 6 4 1 4 11 0 4 6 8 -1 9 3 8 9 7 2 4 4 1 4 1 46 6 2 6 3 8 14 4 6 5 7 3 8 5 8 18 10 3 2 10 2 5 6 3 7 -1 3 12 5 1 14 1 6 9 14 1 6 2 5 6 -3 5 59 -1 8 6 -3 8 1 0 10 0 6 8 9 0 0 2 7 -1 10 -1 4 1 2 4 4 3 4 4 4 9 3 4 1 4 6 14 2 9 3 3 -1 3 6 8 5 8 3 4 3 5 5 4 3 1 60 8 4 2 1 11 2 4 3 6 1 1 6 6 4 3 1 2 12 12 1 7 4 0 6 1 5 3 1 0 9 -1 1 6 6 3 3 5 3 3 6 7 1 3 9 4 1 4 8 6 5 4 3 8 3 1 -4 5 0 3 4 1 8 -1 3 5 2 5 -4 1 3 2 4 2 10 4 2 6 1 6 3 7 2 6 2 12 4 2 12 4 2 4 4 3 0 7 3 6 4 4 -2 15 9 6 5 3 6 3 1 -3 3 25 1 -1 5 4 5 6 2 -1 3 3 14 11 3 0 2 1 3 8 1 5 7 17 -1 5 4 1 1 57 0 -1 -4 -1 16 6 0 3 3 0 5 8 2 6 36 1 -6 28 1 5 1 5 1 0 7 6 8 3 9 3 3 9 15 2 2 13 2 1 5 3 2 0 2 0 -1 2 3 7 4 12 6 5 8 2 1 11 4 1 3 8 6 0 10 -2 7 0 9 4 1 4 2 4 8 16 0 0 4 4 9 15 2 3 -1 7 7 -2 8 5 3 1 5 5 0 7 8 2 2 4 9 6 4 9 3 11 5 -1 3 7 5 8 1 2 0 2 0 5 4 7 2 1 6 2 2 5 4 3 4 2 2 5 2 0 2 2 4 1 5 3 2 1 0 11 2 5 4 18 11 7 6 5 7 4 3 6 2 4 1 -1 0 1 2 7 5 0 8 4 1 8 2 5 2 4 -4 3 2 5 -1 2 7 7 3 3 15 1 2 3 6 11 0 11 5 8 11 1 3 14 0 5 3 6 6 2 1 9 4 2 2 4 2 3 7 9 8 -1 4 11 5 4 9 3 5 5 6 0 1 0 5 0 4 8 2 13 8 3 2 3 3 27 6 6 6 3 10 7 1 -1 6 5 4"
$ws.Range("C2").Value = "5 8 5 8 1 1 3 5 4 2 0 21 12 9 2 -1 16 14 4 9 2 2 3 8 4 4 7 4 1 5 0 4 1 25 5 8 1 3 4 4 0 2 1 6 0 4 2 10 49 7 2 2 3 9 4 8 4 3 2 8 1 5 1 0 0 4 0 4 4 8 2 0 1 -2 -3 3 2 5 5 0 13 5 3 0 4 6 1 3 4 2 1 4 2 5 2 3 2 4 3 11 2 1 4 2 2 25 -3 1 4 14 5 7 3 5 3 12 7 2 2 1 19 4 -1 3 3 3 4 1 0 3 3 4 3 1 -1 13 6 9 4 1 14 3 2 12 0 2 7 21 3 4 9 7 3 1 37 2 4 2 3 0 1 0 1 3 17 3 3 2 7 -2 6 4 -1 2 5 3 8 3 -1 6 2 11 4 1 -1 8 2 0 3 5 1 6 6 6 6 20 1 3 3 10 -4 2 11 3 3 3 1 0 8 3 1 4 9 1 -3 19 8 5 3 4 3 2 0 1 1 5 8 4 18 7 0 8 5 -1 4 0 14 2 0 -2 6 2 0 1 4 0 5 3 57 -1 5 3 3 11 9 2 0 6 -3 4 4 10 4 8 6 3 3 8 8 -2 2 11 11 2 3 3 3 12 3 -4 5 1 3 3 10 13 2 5 5 1 3 5 0 3 4 11 2 1 3 3 4 3 6 2 6 -1 2 12 3 5 3 3 5 2 2 0 3 9 2 1 3 13 4 6 0 5 1 10 1 7 0 5 9 4 3 1 5 1 5 3 10 -2 3 4 5 45 3 4 3 10 3 6 8 1 9 5 3 3 3 3 2 6 -1 6 8 2 4 2 1 5 -1 3 5 6 4 3 -1 2 7 0 4 5 3 0 8 13 3 8 2 6 1 3 6 10 4 13 6 3 6 5 0 5 7 4 4 3 3 1 3 4 2 0 1 20 1 5 9 2 5 2 3 -1 2 2 4 -1 1 6 1 11 5 15 3 3 6 2 8 2 10 5 3 9 5 9 3 11 1 1 6 6 14 2 -1 3 2 1 1 1 2 7 4 0 2 6 3 1 4 2 8 5 2 6 3 3 8 3 1 12 2 17 3 2 0 26 4 -3 0 9 2 2 0 1 2 2 16 8 4 0 7 9 5 1 1 3 -1 10 0 2 4 5 4 8 -2 2 3 11 6 1 3 4 1 0 7 10 -1 9 4 0 6 5 7 11 1 -1 4 9 8 2 3 -2 10 4 4 1 9 2 0 0 1 25 5 4 1 3 7 3 4 2 13 4 2 2 11 3 6 1 2 14 3 16 3 8 5 4 5 2 2 6 1 4 38 2 6 4 5 3 0 2 6 27 12 5 2 8 7 13 4 2 3 1 11 -1 4 5 4 2 2 2 7 6 9 0 0 9 1 3 7 7 2 6 0 0 6 4 7 7 7 2 4 10 2 -4 -1 -1 2 3 2 2 7 10 6 0 5 4 5 1 4 7 1 12 5 2 4 -1 6 -1 3 3 1 2 5 0 -1 4 0 7 4 3 5 13 2 4 0 7 2 0 4 3 7 12 5 6 6 3 1 2 3 7 5 7 2 0 3 5 3 5 6 0 5 8 2 1 5 1 20 1 1 3 8 3 8 7 0 20 4 5 2 2 3 11 17 2 4 1 7 2 -2 5 3 16 4 6 8 12 6 13 4 -1 6 10 7 3 15 2 8 1 1 1 27 11 17 3 8 3 12 5 8 13 8 15 13 3 3 12 8 30 4 10 4"
$ws.Range("B3").Value = "6 6 7 9 15 -1 -1 14 12 7 10 1 6 3 8 24 45 8 1 12 3 39 1 38 1 4 4 40 6 15 4 22 11 4 10 7 17 8 8 11 30 14 25 41 9 -6 9 5 12 9 10 5 3 48 5 72 16 5 18 14 5 7 3 12 1 23 1 27 8 6 9 29 19 2 17 8 14 15 19 9 3 2 13 8 8 23 21 6 6 5 -3 13 20 20 8 12 12 11 15 16 7 9 5 5 18 4 -1 6 17 9 8 0 3 8 1 28 36 3 31 15 16 6 15 1 6 -1 9 5 45 6 16 8 5 4 1 5 16 14 6 5 1 45 19 6 6 7 15 5 5 2 26 27 7 7 7 17 5 15 36 11 5 4 7 1 2 11 52 49 3 11 1 8 1 24 9 7 3 34 6 5 6 10 1 6 6 78 11 -3 22 12 3 18 11 3 11 22 4 14 31 23 7 18 6 33 8 4 8 25 20 5 23 8 12 5 5 4 3 3 9 14 10 11 1 7 47 14 18 7 8 0 11 2 8 5 12 8 15 11 5 8 8 15 4 2 5 0 39 24 7 9 24 10 12 9 8 6 2 36 14 5 25 42 1 31 9 13 5 11 7 8 3 -1 5 4 15 56 10 3 5 12 14 11 4 8 7 3 14 14 -1 7 21 6 5 15 11 4 7 29 6 9 12 7 1 1 14 2 14 6 12 4 7 14 6 4 -2 10 22 2 2 7 8 13 7 5 -1 7 17 21 6 7 23 5 9 6 8 11 8 5 27 32 15 9 3 15 3 7 72 11 9 13 15 6 7 9 8 3 17 14 6 8 12 27 9 6 8 1 1 14 11 6 12 21 33 21 9 7 4 4 4 58 12 8 14 0 2 23 15 5 50 12 4 5 9 5 4 23 5 2 10 24 10 16 8 2 4 13 7 5 11 7 2 5 29 9 11 4 6 7 14 4 31 10 19 7 8 5 6 2 7 18 19 4 5 6 9 13 1 10 50 22 14 3 9 5 22 11 11 0 15 8 0 42 11 -4 2 12 2 47 14 1 0 10 3 4 15 19 6 1 9 6 7 12 25 17 5 21 7 16 4 8 12 6 6 23 7 19 10 25 26 11 1 28 2 11 5 15 17 -4 34 14 0 12 13 24 -1 3 15 24 5 22 5 7 0 8 8 2 2 59 20 20 15 13 4 0 6 32 1 12 5 11 9 4 10 12 14 13 41 8 15 3 4 12 5 0 12 1 10 3 15 11 7 13 6 19 5 10 2 20 4 10 26 4 25 17 9 7 7 12 4 1 4 15 -2 7 2 19 33 8 5 6 9 12 10 4 7 10 15 5 6 21 35 0 11 -1 14 16 3 8 2 15 5 3 41 9 4 6 12 6 34 3 4 14 8 10 6 4 23 24 4 7 8 5 50 7 37 2 4 13 3 5 39 8 17 12 15 8 25 13 12 17 10 9 13 26 4 18 75 8 15 0 2 7 7 5 1 54 27 43 7 10 6 7 8 28 6 6 5 2 9 2 8 6 7 9 6 6 32 0 7 12 25 12 20 6 13 6 18 38 12 16 8 23 16 12 6 6 8 11 6 4 7 3 17 6 11 31 7 22 10 3 13 3 24 6 1 16 9 7 6 9 25 14 11 23 -6 1 8 4 33 4 6 11 34 1 9 2 -2 10 9 12 1 10 1 17 10 1 9 7 4"
$ws.Range("C3").Value = "7 25 18 12 22 9 -2 5 29 37 10 11 17 22 6 7 9 9 5 12 19 5 29 6 9 8 14 4 13 0 10 12 19 4 9 8 9 17 8 24 36 3 11 -2 21 7 1 8 7 4 7 7 12 8 10 10 13 14 6 18 52 7 3 2 0 18 3 10 9 8 5 13 14 6 1 9 4 10 9 18 16 6 6 11 4 16 4 6 13 11 19 14 21 7 7 8 6 9 10 13 2 6 4 6 0 7 8 25 3 12 31 6 18 5 6 16 5 8 12 8 24 28 14 9 11 8 8 16 3 6 9 30 1 5 -3 15 13 18 8 12 6 8 9 12 13 16 7 12 6 16 12 50 7 6 3 43 4 5 18 3 5 6 1 47 10 9 36 41 8 11 7 9 4 8 7 11 6 8 1 7 5 18 12 2 7 13 21 17 7 10 9 7 8 14 9 13 45 0 9 0 6 5 20 16 17 1 9 9 9 23 6 20 8 8 14 15 3 8 20 16 7 12 12 9 6 13 1 13 10 12 12 6 16 8 4 7 10 9 4 9 3 2 7 11 8 10 41 34 6 20 17 21 5 8 16 15 1 4 16 17 8 23 10 10 10 15 8 16 -1 5 6 12 13 2 4 1 23 5 2 37 15 8 5 16 20 7 -6 6 18 7 9 6 9 15 13 5 14 40 11 12 10 12 12 8 12 9 11 10 19 6 6 7 8 7 5 5 1 6 6 1 0 5 14 10 5 2 11 2 53 3 5 5 14 2 6 4 12 4 6 5 3 1 3 28 19 4 21 5 6 4 4 6 16 4 5 16 8 5 9 5 8 -3 27 14 15 5 -3 14 52 7 19 1 6 39 8 16 -1 6 11 29 9 14 3 17 2 3 3 8 5 6 4 10 8 55 8 46 9 5 6 4 5 11 4 10 8 20 1 6 15 7 2 5 3 5 7 7 5 24 16 7 8 7 6 4 16 8 16 10 11 5 7 19 1 0 16 8 0 16 12 9 3 39 12 12 12 3 9 2 45 6 8 10 4 17 4 6 12 11 11 5 0 4 11 10 7 9 30 3 4 2 7 9 6 7 2 10 4 7 5 70 28 14 5 8 3 10 2 7 19 9 8 8 5 21 15 11 20 3 5 3 12 8 11 20 21 8 10 5 15 7 6 12 40 0 5 5 14 4 17 28 9 12 5 7 10 3 11 2 19 10 3 12 13 6 15 4 8 58 6 9 8 3 6 8 11 20 9 9 4 9 23 29 9 6 22 2 7 10 4 2 1 7 2 7 12 2 25 6 14 10 13 22 9 13 12 -6 6 28 1 9 9 7 5 7 18 8 6 -4 4 8 10 8 9 5 43 3 18 56 15 10 5 7 15 6 7 9 21 5 19 6 7 7 3 23 6 19 26 4 5 79 5 11 22 3 19 6 1 8 9 54 13 7 10 7 7 7 11 14 19 8 3 6 46 6 54 6 2 5 6 13 5 11 4 3 22 5 3 43 5 5 4 22 8 4 8 7 9 2 7 1 9 12 5 5 3 5 5 5 3 3 6 12 3 7 40 12 8 8 11 24 14 1 16 5 3 10 10 12 3 11 10 5 7"

# --- OFF sheet: season totals after Week 16 ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 415
$ws.Range("F2").Value = 122
$ws.Range("G2").Value = 142
$ws.Range("J2").Value = 50
$ws.Range("L2").Value = 591
$ws.Range("M2").Value = 412
$ws.Range("O2").Value = 38
$ws.Range("P2").Value = 24
$ws.Range("Q2").Value = 1082
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 343
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 45
$ws.Range("F3").Value = 213
$ws.Range("G3").Value = 91
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 120
$ws.Range("J3").Value = 121

# --- DEF sheet: season totals after Week 16 ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 368
$ws.Range("D2").Value = 22
$ws.Range("F2").Value = 107
$ws.Range("G2").Value = 113
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 73
$ws.Range("L2").Value = 580
$ws.Range("M2").Value = 344
$ws.Range("O2").Value = 53
$ws.Range("P2").Value = 27
$ws.Range("Q2").Value = 1023
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 348
$ws.Range("D3").Value = 13
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 237
$ws.Range("G3").Value = 71
$ws.Range("H3").Value = 44
$ws.Range("I3").Value = 114
$ws.Range("J3").Value = 96
$ws.Range("N3").Value = 51

# --- ST sheet: season totals + distance/return logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 189
$ws.Range("D2").Value = 99
$ws.Range("F2").Value = 711
$ws.Range("G2").Value = 694
$ws.Range("J2").Value = 255
$ws.Range("K2").Value = 232
$ws.Range("B3").Value = 87
$ws.Range("B4").Value = "68 68 69 66 68 68 70 46 66 67 55 59 63 52 64 58 66 47 46 66 69 60 66 70 54 49 65 72 47 64 64 65 66 66 70 69 65 63 62 56 62 62 58 61 59 64 65 69 56 64 65 62 60 60 59 63 57 57 63 65 64 66 73 70 69 69 64 68 66 62 60 67 65 66 66 67 65 70 68 43 57 63 63 64 44 62 66 60 64 61 65 63 66 62 66 64 62 61 64 63 52 65"
$ws.Range("B5").Value = "26 38 27 26 22 22 30 13 26 20 17 27 13 9 30 18 34 12 20 21 7 4 30 33 10 16 21 20 17 19 24 27 46 46 24 71 23 28 21 26 28 31 23 20 17 0 20 33 10 21 27 38 29 22 18 43 9 33 21 68 16 20 38 27 40 44 22 28 26 26 27 30 27 30 24 25 23 20 32 12 12 24 32 23 15 23 20 28 40 0 42 22 14 18 24 26 34 21 24 28 16 15"
$ws.Range("B6").Value = "18 18 15 19 13 17 20 34 17 19 27 26 5 8 16 30 15 20 22 21 14 8 22 33 23 22 20 0 11 24 20 21 15 23 17 26 30 17 41 33 21 21 16 10 16 23 13 26 19 3 32 14 10 14 19 35 21 20 0 9 21 27"
$ws.Range("D3").Value = "40 45 39 54 38 37 44 40 48 51 38 45 48 55 62 47 57 41 38 53 44 57 48 37 41 59 30 39 40 38 50 42 41 54 52 47 57 44 36 57 36 38 55 53 45 42 37 41 38 55 51 38 47 59 32 46 56 38 51 50 63 36 57 45 36 31 36 50 82 49 53 37 58 48 39 55 60 50 41 49 40 58 23 41 53 61 34 57 45 47 22 47 59 40 45 44 56 49 41"
$ws.Range("D4").Value = "0 0 0 0 14 0 0 0 7 10 0 4 0 9 19 0 14 0 0 0 0 13 18 0 8 91 0 4 5 0 0 0 0 0 73 0 0 0 0 0 0 1 19 11 0 0 0 0 0 7 3 0 0 13 0 0 10 0 0 9 14 0 0 14 2 0 11 27 0 0 0 0 0 14 3 0 0 0 0 0 0 8 0 0 13 0 0 34 0 97 0 0 0 0 4 0 15 0 8"
$ws.Range("D5").Value = "0 0 0 0 0 0 6 0 0 0 0 0 11 0 0 0 0 0 0 0 6 0 0 0 9 0 0 0 0 0 0 0 0 0 7 0 0 0 0 0 12 0 0 0 0 0 2 0 0 5 0 0 0 0 0 7 9 8 0 0 0 0 0 17 2 0 0 0 0 6 9 0 0 1 0 0 8 0 16 0 0 0 8 17 0 0 0 0 15 9 0 0 4 0 16 0 0 0 0 0 0 0 0 13 0 15 0 9"

# --- TURNS sheet: season totals after Week 16 ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 23

